# Update Status column (E) values on Sheet1 as part of a "reduction" pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value  = 12
$ws.Range("E4").Value  = 11
$ws.Range("E5").Value  = 0
$ws.Range("E6").Value  = 1
$ws.Range("E7").Value  = 11
$ws.Range("E8").Value  = 11
$ws.Range("E9").Value  = 12
$ws.Range("E10").Value = 11
$ws.Range("E12").Value = 11
$ws.Range("E13").Value = 11
$ws.Range("E15").Value = 12
$ws.Range("E17").Value = 1
$ws.Range("E21").Value = 11
